# Update the "Binning" row (row 6) values for every condition column (B:S)
# from the old binning spec "[-1,60,80,100]" to the new one "[-1,40,70,100]".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6:S6").Value = "[-1,40,70,100]"

# Reflect the author's new cell selection/scroll position in the sheet view
# (previously P9 was selected with the view scrolled to M1; now C13 is
# selected and the view scrolls so C13 is the top-left visible cell).
$ws.Range("C13").Select()
